$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5075  # was 5051
$ws.Range("F6").Value = 5075  # was 5051
$ws.Range("F7").Value = 97  # was 85
$ws.Range("F9").Value = 502  # was 498
$ws.Range("F11").Value = 1147  # was 1143
$ws.Range("F12").Value = 701  # was 695
$ws.Range("F13").Value = 4903  # was 4869
$ws.Range("F15").Value = 54  # was 52
$ws.Range("F16").Value = 72  # was 68
$ws.Range("F17").Value = 205  # was 201
$ws.Range("F18").Value = 210  # was 209
$ws.Range("F19").Value = 95  # was 94
$ws.Range("F20").Value = 241  # was 240
$ws.Range("F21").Value = 3740  # was 3732
$ws.Range("F23").Value = 36  # was 35
$ws.Range("F24").Value = 3608  # was 3583
$ws.Range("F25").Value = 169  # was 165
$ws.Range("F26").Value = 160  # was 157
$ws.Range("F28").Value = 197  # was 192
$ws.Range("F29").Value = 230  # was 227
$ws.Range("F30").Value = 198  # was 197
$ws.Range("F31").Value = 104  # was 103
$ws.Range("F35").Value = 135  # was 134
$ws.Range("F36").Value = 6372  # was 6315
$ws.Range("F37").Value = 1005  # was 998
$ws.Range("F38").Value = 476  # was 475
$ws.Range("F39").Value = 95  # was 94
$ws.Range("F40").Value = 969  # was 968
$ws.Range("F42").Value = 1302  # was 1295
$ws.Range("F44").Value = 629  # was 623
$ws.Range("F46").Value = 2188  # was 2179
$ws.Range("F47").Value = 311  # was 310
$ws.Range("F48").Value = 87  # was 84
$ws.Range("F49").Value = 754  # was 753
$ws.Range("F50").Value = 899  # was 896

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 130  # was 129
$ws.Range("F9").Value = 75  # was 72
$ws.Range("F15").Value = 139  # was 138
$ws.Range("F23").Value = 795  # was 794

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 216  # was 215

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 216  # was 215
$ws.Range("F10").Value = 5075  # was 5051
$ws.Range("F11").Value = 5075  # was 5051
$ws.Range("F12").Value = 97  # was 85
$ws.Range("F15").Value = 75  # was 72
$ws.Range("F16").Value = 701  # was 695
$ws.Range("F17").Value = 4903  # was 4869
$ws.Range("F19").Value = 54  # was 52
$ws.Range("F20").Value = 72  # was 68
$ws.Range("F21").Value = 205  # was 201
$ws.Range("F22").Value = 95  # was 94
$ws.Range("F23").Value = 241  # was 240
$ws.Range("F24").Value = 3608  # was 3583
$ws.Range("F25").Value = 169  # was 165
$ws.Range("F26").Value = 160  # was 157
$ws.Range("F27").Value = 197  # was 192
$ws.Range("F28").Value = 230  # was 227
$ws.Range("F29").Value = 198  # was 197
$ws.Range("F30").Value = 104  # was 103
$ws.Range("F34").Value = 135  # was 134
$ws.Range("F35").Value = 139  # was 138
$ws.Range("F36").Value = 6373  # was 6315
$ws.Range("F37").Value = 1005  # was 998
$ws.Range("F38").Value = 95  # was 94
$ws.Range("F39").Value = 969  # was 968
$ws.Range("F40").Value = 1302  # was 1295
$ws.Range("F42").Value = 629  # was 623
$ws.Range("F44").Value = 2189  # was 2179
$ws.Range("F45").Value = 311  # was 310
$ws.Range("F47").Value = 87  # was 84
$ws.Range("F48").Value = 754  # was 753
$ws.Range("F49").Value = 899  # was 896
